# "add gr and ud data" - append new department/program rows to Sheet1
# and fix a few existing rows whose English-title column had been
# mistakenly duplicated into the Arabic-title column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Fix up three existing rows: column D (arabic_title) had been
#    filled with the English/code value by mistake - replace with the
#    correct Arabic text.
# ---------------------------------------------------------------
$ws.Range("D8").Value  = "اللغة الإنجليزية والترجمة"
$ws.Range("D11").Value = "الحلقوق"
$ws.Range("D12").Value = "الصحه العامه"

# ---------------------------------------------------------------
# 2) Append the new rows (13-31): code, faculty_code, english_title,
#    arabic_title
# ---------------------------------------------------------------
$data = @(
  @("BUSI", "BU", "Business Administration", "إدارة الأعمال"),
  @("ECOM", "BU", "E-Commerce", "التجارة الإلكترونية"),
  @("EMED", "SA", "Electronic Media", "الإعلام الرقمي"),
  @("ENGL", "SA", "English language & Translation", "اللغة الإنجليزية والترجمة"),
  @("FINA", "BU", "Finance", "إدارة مالية"),
  @("ACCT", "BU", "Accounting", "المحاسبة"),
  @("HCAR", "HS", "Healthcare Administration", "إدارة الرعاية الصحية"),
  @("HINF", "HS", "Health Informatics", "المعلوماتية الصحية"),
  @("ISEC", "CI", "Information Security", "أمن المعلومات"),
  @("ITEC", "CI", "Information Technology", "تقنية المعلومات"),
  @("LAW", "SA", "LAW", "القانون"),
  @("MCYS", "CI", "Cyber Security", "الأمن السيبراني"),
  @("PHLT", "HS", "Public Health", "الصحة العامة"),
  @("DATS", "CI", "Data Science", "علوم البيانات"),
  @("TTEC", "SA", "Translation Technologies", "تقنيات الترجمة"),
  @("HQS", "HS", "Executive Master of Healthcare Quality and Patient Safety", "الماجستير التنفيذي لجودة الرعاية الصحية و سلامة المرضى"),
  @("EMBA", "BU", "Executive MBA", "إدارة الاعمال التنفيذي"),
  @("DMKT", "BU", "Digital Marketing", "التسويق الرقمي"),
  @("DENG", "SA", "English Diploma", "دبلوم اللغة الإنجليزية")
)

$startRow = 13
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Range("A$r").Value = $rowVals[0]
    $ws.Range("B$r").Value = $rowVals[1]
    $ws.Range("C$r").Value = $rowVals[2]
    $ws.Range("D$r").Value = $rowVals[3]
}

# ---------------------------------------------------------------
# 3) Give the new rows their own look: a centered "Heading 4"-based
#    style (bold, navy, 12pt Calibri) - build it once on A13, then
#    fan it out across the whole new block via copy/paste-format.
# ---------------------------------------------------------------
$endRow = $startRow + $data.Count - 1

$anchor = $ws.Range("A13")
$anchor.Style = "Heading 4"
$anchor.Font.Size = 12
$anchor.Font.Color = 6299648
$anchor.HorizontalAlignment = -4108

$anchor.Copy()
$ws.Range("A13:D$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 4) Selection / zoom bookkeeping to match the saved view state.
# ---------------------------------------------------------------
$ws.Range("A13:D$endRow").Select()
$excel.ActiveWindow.Zoom = 85
